# Applies the row 16/17/18 species data rotation described in the diff:
#   - Row 16 takes the species/coordinate data that used to live in row 17
#   - Row 17 takes the species/coordinate data that used to live in row 18
#   - Row 18 takes the species/coordinate data that used to live in row 16
#   - Q/R (easting/northing) values are rounded to whole numbers
#   - The Z and AB ("Starttid"/"Sluttid" = 00:00) cells are cleared for all three rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values for rows 16, 17 and 18.
$row16 = @{
    A = $ws.Range("A16").Value2
    B = $ws.Range("B16").Value2
    D = $ws.Range("D16").Value2
    E = $ws.Range("E16").Value2
    F = $ws.Range("F16").Value2
    G = $ws.Range("G16").Value2
    H = $ws.Range("H16").Value2
    Q = $ws.Range("Q16").Value2
    R = $ws.Range("R16").Value2
}

$row17 = @{
    A = $ws.Range("A17").Value2
    B = $ws.Range("B17").Value2
    D = $ws.Range("D17").Value2
    E = $ws.Range("E17").Value2
    F = $ws.Range("F17").Value2
    G = $ws.Range("G17").Value2
    H = $ws.Range("H17").Value2
    Q = $ws.Range("Q17").Value2
    R = $ws.Range("R17").Value2
}

$row18 = @{
    A = $ws.Range("A18").Value2
    B = $ws.Range("B18").Value2
    D = $ws.Range("D18").Value2
    E = $ws.Range("E18").Value2
    F = $ws.Range("F18").Value2
    G = $ws.Range("G18").Value2
    H = $ws.Range("H18").Value2
    Q = $ws.Range("Q18").Value2
    R = $ws.Range("R18").Value2
}

# Write row 17's data into row 16 (rounding the coordinates).
$ws.Range("A16").Value = $row17.A
$ws.Range("B16").Value = $row17.B
$ws.Range("D16").Value = $row17.D
$ws.Range("E16").Value = $row17.E
$ws.Range("F16").Value = $row17.F
$ws.Range("G16").Value = $row17.G
$ws.Range("H16").Value = $row17.H
$ws.Range("Q16").Value = [Math]::Round([double]$row17.Q)
$ws.Range("R16").Value = [Math]::Round([double]$row17.R)

# Write row 18's data into row 17 (rounding the coordinates).
$ws.Range("A17").Value = $row18.A
$ws.Range("B17").Value = $row18.B
$ws.Range("D17").Value = $row18.D
$ws.Range("E17").Value = $row18.E
$ws.Range("F17").Value = $row18.F
$ws.Range("G17").Value = $row18.G
$ws.Range("H17").Value = $row18.H
$ws.Range("Q17").Value = [Math]::Round([double]$row18.Q)
$ws.Range("R17").Value = [Math]::Round([double]$row18.R)

# Write (original) row 16's data into row 18 (rounding the coordinates).
$ws.Range("A18").Value = $row16.A
$ws.Range("B18").Value = $row16.B
$ws.Range("D18").Value = $row16.D
$ws.Range("E18").Value = $row16.E
$ws.Range("F18").Value = $row16.F
$ws.Range("G18").Value = $row16.G
$ws.Range("H18").Value = $row16.H
$ws.Range("Q18").Value = [Math]::Round([double]$row16.Q)
$ws.Range("R18").Value = [Math]::Round([double]$row16.R)

# Clear the now-unused "Starttid"/"Sluttid" (00:00) cells for all three rows.
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
